$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.771.53"
$ws.Range("E2").Value = "  +2.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.559.66"
$ws.Range("E3").Value = "  +5.83%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.02"
$ws.Range("E5").Value = "  +2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.37"
$ws.Range("E6").Value = "  +8.29%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.556.82"
$ws.Range("E9").Value = "  +5.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("E12").Value = "  +1.84%  "

$ws.Range("E13").Value = "  +3.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.25"
$ws.Range("E14").Value = "  +9.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.017.30"
$ws.Range("E15").Value = "  +5.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.656.31"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("E17").Value = "  +3.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.556.27"
$ws.Range("E18").Value = "  +5.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  +4.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.30"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +3.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.92"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.18"
$ws.Range("E24").Value = "  +1.86%  "

$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("E26").Value = "  +6.09%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("E29").Value = "  +4.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  +13.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0843"
$ws.Range("E31").Value = "  +7.82%  "

$ws.Range("E32").Value = "  +4.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.36"
$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  +13.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "424.34"
$ws.Range("E35").Value = "  +12.58%  "

$ws.Range("E36").Value = "  +3.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.16"
$ws.Range("E37").Value = "  +3.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.49"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +5.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.49"
$ws.Range("E42").Value = "  +3.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "156.76"
$ws.Range("E43").Value = "  +7.76%  "

$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.10"
$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.613"
$ws.Range("E46").Value = "  +4.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0535"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.93"
$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  +6.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.89"
$ws.Range("E51").Value = "  +12.04%  "
